$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells are stored as text in the source data (e.g. "259.71"),
# not numbers, so force text format before writing the new value to avoid
# Excel auto-converting the numeric-looking string into a Number cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "259.47"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.177"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06095"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.722"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.483"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.361"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.7992"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08062"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03320"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03045"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09306"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.918"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001688"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006144"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006188"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001099"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.003389"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001502"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.695"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.260"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04588"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007177"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1117"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003133"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01064"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005940"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7505"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06420"
